# Word2016-Tp1-Sobre.docx edit
# Replaces the MERGEFIELD-driven "Remite" (return address) block with
# literal "UNLaM" text, and turns the former literal "Direccion" (address)
# block into the MERGEFIELD-driven one (APELLIDO/NOMBRE, COD_POST/DIRECCION,
# LOCALIDAD), matching an Excel-based mail merge envelope layout.

function New-FlatOpc([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: paragraph 1 ("APELLIDO"/"NOMBRE" merge fields, Remitedesobre)
#         becomes the literal "  UNLaM " text.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$body1 = '<w:body>' + `
  '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Remitedesobre"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>UNLaM</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '</w:p>' + `
'</w:body>'
$p1.Range.InsertXML((New-FlatOpc $body1))

# ------------------------------------------------------------------
# Step 2: paragraph 2 ("COD_POST"/"DIRECCIÓN" merge fields, Remitedesobre)
#         becomes the literal "  Florencio Varela 1903 " text.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$body2 = '<w:body>' + `
  '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Remitedesobre"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">  Florencio Varela 1903 </w:t></w:r>' + `
  '</w:p>' + `
'</w:body>'
$p2.Range.InsertXML((New-FlatOpc $body2))

# ------------------------------------------------------------------
# Step 3: paragraph 3 ("LOCALIDAD" merge field, Remitedesobre) becomes
#         the literal "  (1754) San Justo. " text, and picks up the
#         _GoBack bookmark (previously on its own trailing paragraph).
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$body3 = '<w:body>' + `
  '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Remitedesobre"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">  (1754) San Justo. </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
'</w:body>'
$p3.Range.InsertXML((New-FlatOpc $body3))

# ------------------------------------------------------------------
# Step 4: insert three brand-new "Remitedesobre" paragraphs, framed
#         (the print-frame used for the envelope address block), right
#         after paragraph 3: APELLIDO/NOMBRE, COD_POST/DIRECCIÓN and
#         LOCALIDAD merge fields (showing the "«FIELD»" placeholders,
#         as Word renders an un-merged MERGEFIELD).
# ------------------------------------------------------------------
$rPr24 = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPr24NoProof = '<w:rPr><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$framePr = '<w:framePr w:w="7920" w:h="1980" w:hRule="exact" w:hSpace="141" w:wrap="auto" w:hAnchor="page" w:xAlign="center" w:yAlign="bottom"/>'
$ind = '<w:ind w:left="2126" w:firstLine="851"/>'
$pPrCommon = '<w:pPr><w:pStyle w:val="Remitedesobre"/>' + $framePr + $ind + $rPr24 + '</w:pPr>'

function MergeFieldRuns([string]$fieldName) {
    return '<w:r>' + $rPr24 + '<w:fldChar w:fldCharType="begin"/></w:r>' + `
           '<w:r>' + $rPr24 + '<w:instrText xml:space="preserve"> MERGEFIELD "' + $fieldName + '" </w:instrText></w:r>' + `
           '<w:r>' + $rPr24 + '<w:fldChar w:fldCharType="separate"/></w:r>' + `
           '<w:r>' + $rPr24NoProof + '<w:t>«' + $fieldName + '»</w:t></w:r>' + `
           '<w:r>' + $rPr24 + '<w:fldChar w:fldCharType="end"/></w:r>'
}

$paraApellidoNombre = '<w:p>' + $pPrCommon + `
    (MergeFieldRuns 'APELLIDO') + `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    (MergeFieldRuns 'NOMBRE') + `
  '</w:p>'

$paraCodPostDireccion = '<w:p>' + $pPrCommon + `
    (MergeFieldRuns 'COD_POST') + `
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    (MergeFieldRuns 'DIRECCIÓN') + `
  '</w:p>'

$paraLocalidad = '<w:p>' + $pPrCommon + (MergeFieldRuns 'LOCALIDAD') + '</w:p>'

# A zero-length range right after paragraph 3's paragraph mark is the
# insertion point; a trailing empty <w:p/> keeps the new block's last
# paragraph mark from merging into paragraph 4 ("UNLaM").
$insertPoint = $d.Range($p3.Range.End, $p3.Range.End)
$body4 = '<w:body>' + $paraApellidoNombre + $paraCodPostDireccion + $paraLocalidad + '<w:p/>' + '</w:body>'
$insertPoint.InsertXML((New-FlatOpc $body4))

# ------------------------------------------------------------------
# Step 5: clean up. After the insert above, paragraphs are:
#   1 UNLaM, 2 Florencio Varela, 3 San Justo (+bookmark),
#   4 APELLIDO/NOMBRE, 5 COD_POST/DIRECCIÓN, 6 LOCALIDAD,
#   7 (blank, from trailing <w:p/>), 8 "  UNLaM " (Direccinsobre),
#   9 "  Florencio Varela 1903 " (Direccinsobre),
#   10 "  (1754) San Justo." (Direccinsobre), 11 bookmark-only blank.
# Paragraph 7 (the throwaway blank) needs to disappear, paragraph 8's
# runs need to be cleared (its pPr survives as the final, empty
# paragraph), and paragraphs 9-11 need to be removed outright.
# ------------------------------------------------------------------

# Drop the throwaway blank paragraph left over from the insert.
$d.Paragraphs.Item(7).Range.Delete()

# Remove the last two literal-address paragraphs (Florencio Varela,
# San Justo) and the trailing bookmark-only paragraph completely.
$rngRemove = $d.Range($d.Paragraphs.Item(9).Range.Start, $d.Paragraphs.Item(11).Range.End)
$rngRemove.Delete()

# Clear the remaining "UNLaM" (Direccinsobre) paragraph's runs so only
# its empty paragraph mark (and pPr) survives.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = ""

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($pp.Range.Text)]"
}

# ------------------------------------------------------------------
# Step 6: recipientData.xml gains two more <wne:recipientData><wne:active
# wne:val="1"/></wne:recipientData> entries (mail-merge recipient list
# grew once the Excel data source got wired in).
# ------------------------------------------------------------------
try {
    $mm = $d.MailMerge
    $mm.DataSource.ActiveRecord = 1
} catch {
}
